# MyHotelApp_UserStories.docx - rewording of several user-story paragraphs
# plus relocation of the (hidden) "_GoBack" bookmark from the first story
# paragraph to the last story paragraph in this block.
#
# "navbar switched to partial view" (commit message is not literally about
# this content; we only apply the textual diff shown for the document).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "make an account" paragraph: drop the "(if I plan on making regular
#    reservations)" aside and append a new sentence about account benefits.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "hotel application (if I plan on making regular reservations) or make a reservation without an account.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "hotel application or make a reservation without an account. Creating an account is for the purpose of increasing the brevity of future reservations.",
    2) | Out-Null

# ---------------------------------------------------------------------
# 2) confirmation-email paragraph: add the SMTP-server parenthetical.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "my reservation is successfully completed.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "my reservation is successfully completed (using google SMTP mail server).",
    2) | Out-Null

# ---------------------------------------------------------------------
# 3) text-confirmation paragraph: add the twilio parenthetical.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "my room is ready to be checked into.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "my room is ready to be checked into (using twilio).",
    2) | Out-Null

# ---------------------------------------------------------------------
# 4) pay-online paragraph: add the stripe parenthetical.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "ty to pay for my reservation online.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "ty to pay for my reservation online (using stripe).",
    2) | Out-Null

# ---------------------------------------------------------------------
# 5) "view status of all rooms" paragraph: mention vacancy/cleanliness.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "I want the ability to view the status of all rooms in the hotel.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "I want the ability to view the vacancy status and cleanliness status of all rooms in the hotel.",
    2) | Out-Null

# ---------------------------------------------------------------------
# 6) "rooms available for a given set of dates" paragraph: reword and add
#    the querying/sorting parenthetical.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "available for a given set of dates.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "available for a specified set of dates (by querying and sorting the reservations table by date to populate a list of available rooms).",
    2) | Out-Null

# ---------------------------------------------------------------------
# 7) Relocate the hidden "_GoBack" bookmark from the end of the first
#    story ("...available.") to the end of the last story in this block
#    ("...change room availability to ready.").
# ---------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$found = $d.Content
$found.Find.Execute("change room availability to ready.") | Out-Null
$endRange = $d.Range($found.End - 1, $found.End)
$d.Bookmarks.Add("_GoBack", $endRange) | Out-Null
